$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("B2", 8.40942469049587),
    @("C2", 6.130479764429754),
    @("D2", 5.983608504630558),
    @("E2", 12.93255178033751),
    @("G2", 29.31653884102824),
    @("H2", 14.4120164427717),
    @("K2", 8.094003280954288),
    @("M2", 13.22193071412816),
    @("O2", 22.01687124542553),
    @("B3", 8.107977518448136),
    @("C3", 6.041535734697083),
    @("D3", 5.864388629054091),
    @("E3", 12.72634511710361),
    @("G3", 29.33387441203656),
    @("H3", 14.45597623849262),
    @("K3", 7.823340099493381),
    @("M3", 13.04565019629596),
    @("O3", 22.0781188278254),
    @("B4", 7.918338248044209),
    @("C4", 5.98577116089267),
    @("D4", 5.791695349216153),
    @("E4", 12.602625625596),
    @("G4", 29.35433509050244),
    @("H4", 14.48538446662642),
    @("K4", 7.650654959413193),
    @("M4", 12.93938290011327),
    @("O4", 22.12070061569829),
    @("B5", 7.840037135624149),
    @("C5", 5.962772640433886),
    @("D5", 5.762246205953688),
    @("E5", 12.5530072189384),
    @("G5", 29.36513438328576),
    @("H5", 14.49797578978954),
    @("K5", 7.578712660653009),
    @("M5", 12.8966257288588),
    @("O5", 22.13930086906073),
    @("B6", 7.826977402682793),
    @("C6", 5.958937675390703),
    @("D6", 5.757368063289667),
    @("E6", 12.54481837844238),
    @("G6", 29.36707601835174),
    @("H6", 14.50010323054722),
    @("K6", 7.56667371274783),
    @("M6", 12.88956051741537),
    @("O6", 22.14246469252913),
    @("B7", 7.917286217777505),
    @("C7", 5.985462083009321),
    @("D7", 5.791297422049769),
    @("E7", 12.60195312952354),
    @("G7", 29.35447077814308),
    @("H7", 14.48555181960901),
    @("K7", 7.649690997370263),
    @("M7", 12.93880397627731),
    @("O7", 22.12094641733908),
    @("B8", 8.306499799082662),
    @("C8", 6.10006042927118),
    @("D8", 5.942423528107053),
    @("E8", 12.86089448901182),
    @("G8", 29.32047502031627),
    @("H8", 14.42667190488907),
    @("K8", 8.002062174605323),
    @("M8", 13.16077205764083),
    @("O8", 22.03695480652242),
    @("B9", 9.028857565415144),
    @("C9", 6.314972723505491),
    @("D9", 6.240822830639241),
    @("E9", 13.38844718640194),
    @("G9", 29.33193152342256),
    @("H9", 14.33040526618568),
    @("K9", 8.639119658331182),
    @("M9", 13.60934366737466),
    @("O9", 21.91186709485767),
    @("B10", 9.529024662037374),
    @("C10", 6.466027955921532),
    @("D10", 6.458718510516329),
    @("E10", 13.78380974813643),
    @("G10", 29.38816759945226),
    @("H10", 14.27140918490796),
    @("K10", 9.071510479865076),
    @("M10", 13.94386572681897),
    @("O10", 21.84430074741198),
    @("B11", 9.748954743540938),
    @("C11", 6.533091294776041),
    @("D11", 6.557039896666629),
    @("E11", 13.96446084406061),
    @("G11", 29.42413375439297),
    @("H11", 14.24712344933027),
    @("K11", 9.259999734859059),
    @("M11", 14.09646205893879),
    @("O11", 21.8188820786194),
    @("B12", 9.831075584486069),
    @("C12", 6.558235351790062),
    @("D12", 6.594116438501825),
    @("E12", 14.03291086150444),
    @("G12", 29.43924289144418),
    @("H12", 14.2382945127932),
    @("K12", 9.330162330206546),
    @("M12", 14.15425255730278),
    @("O12", 21.81002381617806),
    @("B13", 9.813442078308942),
    @("C13", 6.552831521088271),
    @("D13", 6.586138884619823),
    @("E13", 14.01816823107194),
    @("G13", 29.43592270212773),
    @("H13", 14.24017962645967),
    @("K13", 9.3151060555545),
    @("M13", 14.14180695116892),
    @("O13", 21.81189744025625),
    @("B14", 9.755734526540547),
    @("C14", 6.535165025811502),
    @("D14", 6.560093549213675),
    @("E14", 13.97009185183821),
    @("G14", 29.42534696831452),
    @("H14", 14.24638971797481),
    @("K14", 9.265796566925575),
    @("M14", 14.10121671591541),
    @("O14", 21.81813791367395),
    @("B15", 9.7202337643633),
    @("C15", 6.524310654146396),
    @("D15", 6.544118597298075),
    @("E15", 13.94064692678996),
    @("G15", 29.41906286195092),
    @("H15", 14.2502414605785),
    @("K15", 9.235434024984626),
    @("M15", 14.07635314461721),
    @("O15", 21.82206037251825),
    @("B16", 9.514492690673649),
    @("C16", 6.461610766747862),
    @("D16", 6.452273353325094),
    @("E16", 13.77201351372422),
    @("G16", 29.38602581434398),
    @("H16", 14.27304771581502),
    @("K16", 9.059024101915533),
    @("M16", 13.93389693852038),
    @("O16", 21.84606914652841),
    @("B17", 9.386278940461656),
    @("C17", 6.422713190690382),
    @("D17", 6.395694706067665),
    @("E17", 13.66871397068029),
    @("G17", 29.36841631909006),
    @("H17", 14.28769263599566),
    @("K17", 8.948674585130034),
    @("M17", 13.84657293329335),
    @("O17", 21.8621616268283),
    @("B18", 9.311822129937246),
    @("C18", 6.400185579147223),
    @("D18", 6.363079107860285),
    @("E18", 13.60937825413377),
    @("G18", 29.35926555684508),
    @("H18", 14.29635616754682),
    @("K18", 8.884433663204886),
    @("M18", 13.79638927508838),
    @("O18", 21.87191797562862),
    @("B19", 9.286492362697093),
    @("C19", 6.392531967606089),
    @("D19", 6.352024751031806),
    @("E19", 13.58930412019871),
    @("G19", 29.35633527374754),
    @("H19", 14.29933072012044),
    @("K19", 8.862551514293154),
    @("M19", 13.77940699502687),
    @("O19", 21.87530716594117),
    @("B20", 9.400001696154892),
    @("C20", 6.426870022320228),
    @("D20", 6.401725461557102),
    @("E20", 13.67970270016418),
    @("G20", 29.37018970601184),
    @("H20", 14.28610879916663),
    @("K20", 8.960501514524019),
    @("M20", 13.85586469485616),
    @("O20", 21.86039674794263),
    @("B21", 9.772716675798231),
    @("C21", 6.540361032256035),
    @("D21", 6.567748231067817),
    @("E21", 13.98421250130875),
    @("G21", 29.42841293198831),
    @("H21", 14.24455568392536),
    @("K21", 9.280313160696659),
    @("M21", 14.11313933141355),
    @("O21", 21.81628409357895),
    @("B22", 10.00950885906146),
    @("C22", 6.613062859472309),
    @("D22", 6.675329397150128),
    @("E22", 14.18343601210078),
    @("G22", 29.47514412412587),
    @("H22", 14.21954090979076),
    @("K22", 9.482238581273151),
    @("M22", 14.28129123428625),
    @("O22", 21.79192697104613),
    @("B23", 9.883771142981622),
    @("C23", 6.574399504337778),
    @("D23", 6.618008671927589),
    @("E23", 14.07711123626487),
    @("G23", 29.44941044615442),
    @("H23", 14.23269552956874),
    @("K23", 9.375126015329524),
    @("M23", 14.19156244087474),
    @("O23", 21.80451675887177),
    @("B24", 9.393799953522761),
    @("C24", 6.424991231263194),
    @("D24", 6.398999227323131),
    @("E24", 13.67473452621382),
    @("G24", 29.36938492612555),
    @("H24", 14.2868240918622),
    @("K24", 8.955157047471255),
    @("M24", 13.85166382271691),
    @("O24", 21.86119307815438),
    @("B25", 8.838440866665719),
    @("C25", 6.25797470544139),
    @("D25", 6.160151936886462),
    @("E25", 13.24405716829999),
    @("G25", 29.32044081389635),
    @("H25", 14.35438983559413),
    @("K25", 8.472861898051974),
    @("M25", 13.48689545973703),
    @("O25", 21.94144539941761)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
